# Add a "Span (km)" data column into the output sheet.
# A new column is inserted before the existing "Compensation (%)" column
# (column D), pushing it and everything after it one column to the right,
# and the new column is filled with the span length (25) used for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting Compensation/fiberAeff/fiberAlphadB right.
$ws.Columns("D:D").Insert()

# Header for the new column.
$ws.Range("D1").Value = "Span (km)"

# Fill the new column with the span value (25 km) for every data row (2-25).
$ws.Range("D2:D25").Value = 25
